$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A3").Value = "1 yr"
$ws.Range("B3").Value = "98.5 (98.3,98.6)"
$ws.Range("C3").Value = "97.3 (97.1,97.4)"
$ws.Range("D3").Value = "96.6 (96.4,96.9)"
$ws.Range("F3").Value = "-1.8 (-2.1,-1.5)"
$ws.Range("G3").Value = "0.6 (0.3,0.9)"

$ws.Range("A4").Value = "3 yrs"
$ws.Range("D4").Value = "93.2 (92.8,93.5)"

$ws.Range("A5").Value = "5 yrs"
$ws.Range("B5").Value = "94.9 (94.6,95.3)"

$ws.Range("A7").Value = "1 yr"
$ws.Range("B7").Value = "0.992 (0.991,0.993)"
$ws.Range("C7").Value = "0.985 (0.984,0.986)"
$ws.Range("D7").Value = "0.981 (0.979,0.982)"

$ws.Range("A8").Value = "3 yrs"
$ws.Range("B8").Value = "2.939 (2.933,2.945)"
$ws.Range("C8").Value = "2.897 (2.893,2.902)"
$ws.Range("D8").Value = "2.876 (2.868,2.883)"
$ws.Range("F8").Value = "-0.064 (-0.074,-0.054)"

$ws.Range("A9").Value = "5 yrs"
$ws.Range("B9").Value = "4.852 (4.840,4.864)"
$ws.Range("C9").Value = "4.761 (4.751,4.770)"
$ws.Range("D9").Value = "4.715 (4.700,4.730)"
$ws.Range("E9").Value = "-0.091 (-0.106,-0.076)"
$ws.Range("F9").Value = "-0.137 (-0.157,-0.116)"
